$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B9: "(16-00)20" -> "2 hour 20 min", drop the red-text style (back to default)
$ws.Range("B9").Value = "2 hour 20 min"
$ws.Range("B9").Style = "Normal"

# B10: new note "18-15" with the red-text style (same style B9 used to have)
$ws.Range("B10").Value = "18-15"
$ws.Range("B10").Font.Color = 255

# Column B sized to fit, selection moved to B10, no frozen/top-left override
$ws.Range("B1").ColumnWidth = 10.83
$ws.Range("B10").Select()
